$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regen save_data column G ("K") to use K (strikes/putts?) instead of the old
# "Strike#" values. New values calculated/regenerated for each row.
$newK = @{
    2  = 4
    3  = 6
    4  = 4
    5  = 9
    6  = 7
    7  = 2
    8  = 1
    9  = 6
    10 = 6
    11 = 4
    12 = 2
    13 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
